# "Add files via upload" — the re-uploaded deck drops the extra picture
# and the "Таблица сравнения" caption textbox that had been placed on
# slide 2 ("p19" in the original Google Slides export), while leaving
# the title/body placeholders on that slide untouched.
#
# (The shape-id renumbering visible in the raw OOXML diff for every
# slide/notes-slide after this point is just Google Slides compacting
# its internal id counter on export after these two shapes were
# removed — ids are read-only/host-assigned in the PowerPoint object
# model, so there is nothing to "author" there; deleting the two
# shapes is the actual, reproducible edit.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Walk backwards so deleting a shape doesn't shift the index of the
# shape we still need to inspect/delete.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    $name = $shape.Name

    $isExtraPicture = ($shape.Type -eq 13) -and ($name -like "Google Shape;151;*")
    $isCaptionBox = $false
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Таблица сравнения") {
            $isCaptionBox = $true
        }
    }

    if ($isExtraPicture -or $isCaptionBox) {
        $shape.Delete()
    }
}
